$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 4 (Banca / 180003 / 100) moves down to row 6, and two brand
# new rows (4 and 5) are inserted above it. Clear the old row 4 first so
# its former contents don't linger in cells that aren't re-used below.
$ws.Range("A4:H4").Value = $null

# --- Header row (row 1): insert a new "Ref" column before Dare/Avere ---
$ws.Range("F1").Value = "Ref"
$ws.Range("G1").Value = "Dare"
$ws.Range("H1").Value = "Avere"

# --- Row 2: Prima Alpha S.p.A. (customer) ---
$ws.Range("B2").Value = "Prima Alpha S.p.A."
$ws.Range("C2").Value = 1
$ws.Range("E2").Value = "IT00115719999"
$ws.Range("F2").Value = $null
$ws.Range("G2").Value = 1000

# --- Row 3: Notaio Decimo Jackson (supplier) ---
$ws.Range("B3").Value = "Notaio Decimo Jackson"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = "IT10242670015"
$ws.Range("G3").Value = $null
$ws.Range("H3").Value = 500

# --- Row 4: Latte Beta Due s.n.c. (new customer row, with Ref = RiBA) ---
$ws.Range("A4").Value = 152220
$ws.Range("B4").Value = "Latte Beta Due s.n.c."
$ws.Range("C4").Value = 1
$ws.Range("E4").Value = "IT02345670018"
$ws.Range("F4").Value = "RiBA"
$ws.Range("G4").Value = 150

# --- Row 5: Freie Universität Berlin (new supplier row) ---
$ws.Range("B5").Value = "Freie Universität Berlin"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = "DE123456788"
$ws.Range("H5").Value = 200

# --- Row 6: Banca (moved down from the old row 4) ---
$ws.Range("A6").Value = 180003
$ws.Range("B6").Value = "Banca"
$ws.Range("G6").Value = 100

# Move the selection to A6 to match the saved cursor position
$ws.Range("A6").Select()
